# Insert a new row before row 9 and populate it with the "Evolution of Type-1 Clones" entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Evolution of Type-1 Clones"
$ws.Range("B9").Value = "N. G" + [char]0x00F6 + "de"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "2009"
$ws.Range("D9").Value = "10.1109/SCAM.2009.17"
$ws.Range("E9").Value = "978-0-7695-3793-1"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "https://ieeexplore.ieee.org/stamp/stamp.jsp?arnumber=5279977"
